# Generate Report for Handback
# This script swaps the rows for the two files that have been processed
# (cb3993ee-... stays "Ready for handoff", d9a56058-... is now handed back)
# across the Overview / zh-cn / de-de sheets, and fills in the new
# handback metadata (Latest Target File / Latest Handback File / Latest
# Handback DateTime) for the d9a56058-... row on the locale sheets.

$wb = $excel.ActiveWorkbook

$cbName     = "cb3993ee-91f8-4746-8fb8-551e3e2180c6.md"
$cbPath     = "e2e\cb3993ee-91f8-4746-8fb8-551e3e2180c6.md"
$cbUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c68f37156d17e0a799b57dc69eebfd4f6d5cab4d/e2e/cb3993ee-91f8-4746-8fb8-551e3e2180c6.md"

$d9aName    = "d9a56058-d8ce-4a43-ba81-5082fe05ad0c.md"
$d9aPath    = "e2e\d9a56058-d8ce-4a43-ba81-5082fe05ad0c.md"
$d9aUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c1fc0f93a35a4c305ca9a5431d2c4666752de6a9/e2e/d9a56058-d8ce-4a43-ba81-5082fe05ad0c.md"

# ---------------------------------------------------------------------
# Overview sheet: row 2 becomes d9a (now handed back), row 3 becomes cb
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $d9aName
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G2").Value = "2016-08-26 00:39:30"

$wsOverview.Range("A3").Value = $cbName
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-26 00:39:14"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $d9aUrl, "", "", $d9aPath)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $cbUrl, "", "", $cbPath)

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet: row 2 becomes d9a (handed back), row 3 becomes cb
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $d9aName
$wsZhCn.Range("G2").Value = "d9a56058-d8ce-4a43-ba81-5082fe05ad0c.9c85f379e7b23e14c15d344f4b1879161b4a8256.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-26 00:39:25"
$wsZhCn.Range("J2").Value = "d9a56058-d8ce-4a43-ba81-5082fe05ad0c.9c85f379e7b23e14c15d344f4b1879161b4a8256.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-26 00:39:41"

$wsZhCn.Range("A3").Value = $cbName
$wsZhCn.Range("G3").Value = "cb3993ee-91f8-4746-8fb8-551e3e2180c6.43920cd83f5773bd3ee75461e2f210513a0ce339.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-26 00:39:09"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $d9aUrl, "", "", $d9aName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $d9aUrl, "", "", $d9aName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $cbUrl, "", "", $cbName)

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet: row 2 becomes d9a (handed back), row 3 becomes cb
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $d9aName
$wsDeDe.Range("G2").Value = "d9a56058-d8ce-4a43-ba81-5082fe05ad0c.9c85f379e7b23e14c15d344f4b1879161b4a8256.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-26 00:39:30"
$wsDeDe.Range("J2").Value = "d9a56058-d8ce-4a43-ba81-5082fe05ad0c.9c85f379e7b23e14c15d344f4b1879161b4a8256.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-26 00:39:49"

$wsDeDe.Range("A3").Value = $cbName
$wsDeDe.Range("G3").Value = "cb3993ee-91f8-4746-8fb8-551e3e2180c6.43920cd83f5773bd3ee75461e2f210513a0ce339.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-26 00:39:14"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $d9aUrl, "", "", $d9aName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $d9aUrl, "", "", $d9aName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $cbUrl, "", "", $cbName)

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
